$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the first column header from "modalidade" to "id_processo"
# (row values below it are untouched; this is the only data change)
$ws.Range("A1").Value = "id_processo"

# Reflect the cell that was selected/active when the workbook was last saved
$ws.Range("F24").Select()
